$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","G","H","J","K","L","M")

$values = @(1.462840114874638, 0.03909597882196891, 0.5212896238752194, 0.1617850459137635, 4.684645036334302, 3.341293950262951, 0.05605884564244512, 0.8291815869093, 0.4168552195945665, 0.4023072943066772)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $values[$i]
}

$values = @(1.450873262681398, 0.03759667802400557, 0.5194637119972185, 0.1618235631522111, 4.619952377051931, 3.314593868655209, 0.05586085644171401, 0.8170858302977848, 0.4166258814617549, 0.4004843631496975)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $values[$i]
}

$values = @(1.444492841053375, 0.03665206546729749, 0.5185302292796621, 0.1618820773932921, 4.581139980676852, 3.298798990830988, 0.05573350597198257, 0.8103379899158938, 0.4166555634710321, 0.3995931291012411)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $values[$i]
}

$values = @(1.442136229830965, 0.03626102529050357, 0.5181970707447192, 0.1619147057539685, 4.565551681334568, 3.292512954018179, 0.05568014804490318, 0.8077589847301425, 0.416710597049601, 0.3992873581519589)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $values[$i]
}

$values = @(1.441759627906322, 0.03619572304462793, 0.5181446052201579, 0.1619206546444758, 4.562977014959614, 3.291478248487238, 0.05567119951908506, 0.8073410601899838, 0.4167223306163521, 0.3992400548398258)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $values[$i]
}

$values = @(1.444460072888177, 0.03664681655025603, 0.5185255448178623, 0.1618824818454598, 4.580928828720687, 3.29871360596016, 0.0557327922951254, 0.8103025169356641, 0.4166561317319974, 0.3995887728082899)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "7").Value = $values[$i]
}

$values = @(1.458513315577079, 0.03858397284178849, 0.5206211406459573, 0.1617911001292107, 4.662149957451732, 3.331963327346045, 0.05599177523296994, 0.8248700660718953, 0.4167407950604272, 0.4016314539965222)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $values[$i]
}

$values = @(1.493740530991317, 0.04219476107989095, 0.5262174432228193, 0.1618879004839364, 4.828672949118754, 3.401930511998785, 0.05645406687424437, 0.8588257200107137, 0.4182572547152716, 0.4074444984251109)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $values[$i]
}

$values = @(1.524295256288895, 0.04473693970265913, 0.5312339289915116, 0.1621265488784367, 4.955499191467425, 3.456262470215734, 0.05676638214067697, 0.8870648165081434, 0.4201921811055627, 0.4128155415928489)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "10").Value = $values[$i]
}

$values = @(1.539210316648735, 0.04587026369598846, 0.5337123101557353, 0.1622713506626745, 5.01418338459024, 3.481620377039576, 0.05690262483986963, 0.9006282079105574, 0.4212501966082129, 0.4154976392904146)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "11").Value = $values[$i]
}

$values = @(1.545004207498465, 0.04629616416309545, 0.53467900634584, 0.1623313808134483, 5.036548826708668, 3.491315366269646, 0.05695338541384309, 0.9058675242530967, 0.4216763595351125, 0.4165475773418095)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "12").Value = $values[$i]
}

$values = @(1.543749903716673, 0.04620458310802178, 0.5344695581640622, 0.1623182213410175, 5.031725658353139, 3.489223260031451, 0.05694249009072472, 0.9047345556872415, 0.4215834437518282, 0.4163199303140956)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "13").Value = $values[$i]
}

$values = @(1.539684060719992, 0.04590536786479049, 0.5337912761987553, 0.1622761853078956, 5.016020533304697, 3.482416134768357, 0.05690681756931415, 0.9010571821715416, 0.4212847462215876, 0.4155833314708914)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "14").Value = $values[$i]
}

$values = @(1.537212608053352, 0.04572166653109377, 0.5333794782653172, 0.1622511133861408, 5.00641933322072, 3.478258625577951, 0.05688485902913687, 0.8988181188632609, 0.4211051066643279, 0.4151366069635145)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "15").Value = $values[$i]
}

$values = @(1.523340944248616, 0.0446624141503662, 0.5310759079513332, 0.1621178137422774, 4.951684020645928, 3.45461820569551, 0.05675736161575884, 0.8861928512567658, 0.4201266106830559, 0.4126450607386616)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "16").Value = $values[$i]
}

$values = @(1.515091124428011, 0.04400671720284066, 0.5297129987037437, 0.162045310262922, 4.918359749691376, 3.440280137594669, 0.05667765753266174, 0.8786313770564504, 0.4195718377552566, 0.4111776991594525)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "17").Value = $values[$i]
}

$values = @(1.510441644611092, 0.04362739730705556, 0.5289475742145271, 0.1620070200520658, 4.899285685473757, 3.43209367983647, 0.05663126459007284, 0.8743497303320851, 0.41926948625688, 0.4103561899284429)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "18").Value = $values[$i]
}

$values = @(1.508883834609179, 0.04349858978352472, 0.5286915909338745, 0.1619946420762872, 4.892843519517072, 3.429332257684933, 0.05661546218946789, 0.872911633596118, 0.4191699918474825, 0.4100819032728964)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "19").Value = $values[$i]
}

$values = @(1.515959438935369, 0.04407674247960358, 0.5298561699653703, 0.1620526753706066, 4.921897527596997, 3.441800194965339, 0.05668619896449822, 0.8794293214919549, 0.4196291622255472, 0.4113315762595988)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "20").Value = $values[$i]
}

$values = @(1.540874339388694, 0.0459933427547341, 0.5339897394495523, 0.1622883913818107, 5.020629623059619, 3.484413040602817, 0.05691731796726973, 0.90213451622202, 0.4213717888635387, 0.4157987582943861)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "21").Value = $values[$i]
}

$values = @(1.558007908251483, 0.04722695078961436, 0.5368555398563188, 0.1624727301669253, 5.085990622918359, 3.512802245541764, 0.05706352372037493, 0.9175749456651658, 0.4226593956774991, 0.4189181419439905)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "22").Value = $values[$i]
}

$values = @(1.548785652135848, 0.04657026944297371, 0.5353109929922795, 0.1623715787178845, 5.051029726050388, 3.497600991981642, 0.05698593188858236, 0.9092790764284473, 0.4219585869128082, 0.4172350007918695)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "23").Value = $values[$i]
}

$values = @(1.515566582846134, 0.04404509138056767, 0.5297913857626213, 0.1620493350359133, 4.920297835031022, 3.44111280044828, 0.05668233915926724, 0.8790683666484256, 0.4196031941384177, 0.4112619395915189)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "24").Value = $values[$i]
}

$values = @(1.483389854692234, 0.04123766562749509, 0.5245444815273714, 0.1618322416534674, 4.782841574439345, 3.382490143734657, 0.05633383685458337, 0.8490622412420237, 0.4177026993107091, 0.4056785996395078)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "25").Value = $values[$i]
}
